# feat: add 2022-Q4 data
#
# 1. Duplicate the existing "2022-Q3" sheet (so the new sheet inherits all
#    formatting) and place the copy immediately before it; rename the copy
#    to "2022-Q4" and overwrite its fund-holding figures with the new
#    quarter's numbers (fund codes/names stay the same).
# 2. The original "2022-Q3" sheet is left completely untouched, so it keeps
#    representing the (unchanged) 2022-Q3 data - it is simply now the 3rd
#    tab instead of the 2nd.
# 3. Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the
#    top of the data table and shift every other quarter down by one row.

$wb = $excel.ActiveWorkbook

# Helper: force a cell to hold a literal TEXT value (matches the workbook's
# convention of storing these figures as text, not numbers) without leaving
# a stray number-format behind on the cell's style.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Step 1: duplicate "2022-Q3" -> new sheet before it, renamed "2022-Q4"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Update the fund-holding figures on the new 2022-Q4 sheet.
Set-TextValue $q4.Range("D2") "1.98"
Set-TextValue $q4.Range("F2") "2.52"
Set-TextValue $q4.Range("G2") "0.0499"
$q4.Range("H2").Value = 10

Set-TextValue $q4.Range("D3") "0.60"
Set-TextValue $q4.Range("E3") "97.54"
Set-TextValue $q4.Range("F3") "2.51"
$q4.Range("H3").Value = 10

Set-TextValue $q4.Range("D4") "0.27"
Set-TextValue $q4.Range("E4") "98.47"
Set-TextValue $q4.Range("F4") "2.53"
Set-TextValue $q4.Range("G4") "0.0068"
$q4.Range("H4").Value = 10

# ---------------------------------------------------------------------
# Step 2: refresh the "总计" summary table with the new quarter on top.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A8").NumberFormat = $summary.Range("A7").NumberFormat
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

$rows = @(
    @(0, "2022-Q4", 3, 0.07000000000000001),
    @(1, "2022-Q3", 3, 0.08),
    @(2, "2022-Q2", 2, 0.11),
    @(3, "2022-Q1", 3, 0.17),
    @(4, "2021-Q4", 4, 3.85),
    @(5, "2021-Q3", 11, 5.46),
    @(6, "2021-Q2", 5, 1.74)
)

$r = 2
foreach ($row in $rows) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
